$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.139.22"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.565.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.36"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.566.11"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.85"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.173.68"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.08"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.554.29"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.227.18"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.44"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.71"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.26"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.609"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.712.49"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.12"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.87"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.566.16"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.48"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.44"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.10%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.80"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.54"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "173.91"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0844"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.19"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.00"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.02"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.13"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.12"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.93%  "
